$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.250077247619629
$ws.Range("B1").Value = 1.586637377738953
$ws.Range("C1").Value = 3.706434488296509
$ws.Range("D1").Value = 3.382776021957397
$ws.Range("E1").Value = 0.9840773940086365
